# Chris found that the dates weren't being updated between different beta-gal
# data sets, so the "date" column (column K) on Sheet2 and Sheet3 is fixed up
# here to read the correct assay date (2018-07-28) instead of the stale
# 2018-07-24 value that had been carried over from an earlier file.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Column K is "date" on both sheets; data rows run from row 2 through row 41.
$ws2.Range("K2:K41").Value = "2018-07-28"
$ws3.Range("K2:K41").Value = "2018-07-28"
